$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.415.41'
$ws.Range("E2").Value = '  +0.36%  '
# Row 3
$ws.Range("D3").Value = '1.880.32'
$ws.Range("E3").Value = '  +0.42%  '
# Row 4
$ws.Range("E4").Value = '  +0.09%  '
# Row 5
$ws.Range("D5").Value = '''0.7170'
$ws.Range("E5").Value = '  +1.25%  '
# Row 6
$ws.Range("D6").Value = '''243.62'
$ws.Range("E6").Value = '  +0.74%  '
# Row 7
$ws.Range("E7").Value = '  +0.15%  '
# Row 8
$ws.Range("D8").Value = '''0.07936'
$ws.Range("E8").Value = '  +1.74%  '
# Row 9
$ws.Range("D9").Value = '''0.3149'
$ws.Range("E9").Value = '  +1.40%  '
# Row 11
$ws.Range("D11").Value = '''0.08136'
$ws.Range("E11").Value = '  -2.86%  '
# Row 12
$ws.Range("D12").Value = '1.898.69'
$ws.Range("E12").Value = '  +1.01%  '
# Row 13
$ws.Range("D13").Value = '''94.91'
# Row 14
$ws.Range("D14").Value = '''5.242'
$ws.Range("E14").Value = '  +0.16%  '
# Row 15
$ws.Range("D15").Value = '''0.7087'
$ws.Range("E15").Value = '  -1.11%  '
# Row 16
$ws.Range("D16").Value = '''6.396'
$ws.Range("E16").Value = '  +4.28%  '
# Row 17
$ws.Range("D17").Value = '''0.000008430'
$ws.Range("E17").Value = '  +0.56%  '
# Row 18
$ws.Range("D18").Value = '29.436.84'
$ws.Range("E18").Value = '  +0.41%  '
# Row 19
$ws.Range("D19").Value = '''252.91'
$ws.Range("E19").Value = '  +5.26%  '
# Row 20
$ws.Range("D20").Value = '''13.36'
$ws.Range("E20").Value = '  +1.25%  '
# Row 21
$ws.Range("D21").Value = '2.143.83'
$ws.Range("E21").Value = '  +0.79%  '
# Row 22
$ws.Range("D22").Value = '''1.002'
$ws.Range("E22").Value = '  +0.08%  '
# Row 23
$ws.Range("D23").Value = '''7.719'
$ws.Range("E23").Value = '  -0.26%  '
# Row 24
$ws.Range("E24").Value = '  +0.07%  '
# Row 25
$ws.Range("D25").Value = '''0.1586'
$ws.Range("E25").Value = '  -0.47%  '
# Row 26
$ws.Range("D26").Value = '''9.071'
$ws.Range("E26").Value = '  +0.46%  '
# Row 27
$ws.Range("D27").Value = '''162.15'
$ws.Range("E27").Value = '  -0.26%  '
# Row 28
$ws.Range("D28").Value = '''18.94'
$ws.Range("E28").Value = '  +2.56%  '
# Row 29
$ws.Range("D29").Value = '''1.507'
$ws.Range("E29").Value = '  +0.12%  '
# Row 30
$ws.Range("D30").Value = '''4.416'
# Row 31
$ws.Range("D31").Value = '''4.293'
$ws.Range("E31").Value = '  -1.38%  '
# Row 32
$ws.Range("D32").Value = '''1.225'
$ws.Range("E32").Value = '  -0.61%  '
# Row 33
$ws.Range("D33").Value = '''0.05332'
$ws.Range("E33").Value = '  -0.42%  '
# Row 34
$ws.Range("D34").Value = '''1.945'
$ws.Range("E34").Value = '  +0.24%  '
# Row 35
$ws.Range("D35").Value = '''0.7577'
$ws.Range("E35").Value = '  +1.20%  '
# Row 36
$ws.Range("D36").Value = '''1.177'
$ws.Range("E36").Value = '  +0.40%  '
# Row 37
$ws.Range("D37").Value = '''2.699'
$ws.Range("E37").Value = '  +0.53%  '
# Row 38
$ws.Range("D38").Value = '''0.01894'
$ws.Range("E38").Value = '  +0.72%  '
# Row 39
$ws.Range("D39").Value = '1.278.23'
$ws.Range("E39").Value = '  +2.87%  '
# Row 40
$ws.Range("D40").Value = '''2.765'
$ws.Range("E40").Value = '  +1.21%  '
# Row 41
$ws.Range("D41").Value = '''6.405'
$ws.Range("E41").Value = '  -1.92%  '
# Row 42
$ws.Range("D42").Value = '''112.85'
$ws.Range("E42").Value = '  +2.67%  '
# Row 43
$ws.Range("D43").Value = '''0.9083'
$ws.Range("E43").Value = '  +1.77%  '
# Row 44
$ws.Range("D44").Value = '''74.17'
$ws.Range("E44").Value = '  +2.68%  '
# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '''1.002'
$ws.Range("E45").Value = '  +0.12%  '
# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '''0.00000000130'
$ws.Range("E46").Value = '  +0.88%  '
# Row 47
$ws.Range("D47").Value = '2.039.35'
$ws.Range("E47").Value = '  +0.96%  '
# Row 48
$ws.Range("E48").Value = '  +0.91%  '
# Row 49
$ws.Range("D49").Value = '''0.5208'
$ws.Range("E49").Value = '  +0.15%  '
# Row 50
$ws.Range("D50").Value = '''9.508'
$ws.Range("E50").Value = '  +0.54%  '
# Row 51
$ws.Range("E51").Value = '  +0.42%  '
